$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 10, shifting existing rows 10-15 down to 11-16
$ws.Rows.Item(10).Insert()

# Fill the new row 10 with the data for the new weekly price point
$ws.Cells.Item(10, 1).Value = 5
$ws.Cells.Item(10, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(10, 3).Value = "Maule"
$ws.Cells.Item(10, 4).Value = 44447
$ws.Cells.Item(10, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(10, 5).Value = 7
$ws.Cells.Item(10, 6).Value = 100112026
$ws.Cells.Item(10, 7).Value = "Haba"
$ws.Cells.Item(10, 8).Value = "Sin especificar"
$ws.Cells.Item(10, 9).Value = "Primera"
$ws.Cells.Item(10, 10).Value = 200
$ws.Cells.Item(10, 11).Value = 12000
$ws.Cells.Item(10, 12).Value = 12000
$ws.Cells.Item(10, 13).Value = 12000
$ws.Cells.Item(10, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(10, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(10, 16).Value = 480
$ws.Cells.Item(10, 17).Value = 25
$ws.Cells.Item(10, 18).Value = "Hortaliza"
